$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 3..20) {
    $ws.Cells.Item($r, 2).Value = "  "
}

foreach ($r in 22..39) {
    $ws.Cells.Item($r, 2).Value = "  "
}
